$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/date strings) ---
$ws.Range("A8").Value = "Volume 31   Number  45"
$ws.Range("C9").Value = "Report Covering the Week  11/4/2024  Through  11/10/2024"

# --- Data cell updates (rows 14-30) ---
$ws.Range("M14").Value = 0
$ws.Range("D15").NumberFormat = "General"
$ws.Range("D15").Value = "'0"
$ws.Range("E15").NumberFormat = "General"
$ws.Range("E15").Value = "'***.*"
$ws.Range("L15").Value = 0
$ws.Range("C16").Value = 7
$ws.Range("D16").NumberFormat = "#,##0"
$ws.Range("D16").Value = 1
$ws.Range("E16").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("E16").Value = 600
$ws.Range("F16").Value = 11
$ws.Range("G16").Value = 2
$ws.Range("H16").Value = 450
$ws.Range("I16").Value = 107
$ws.Range("J16").Value = 93
$ws.Range("K16").Value = 15.053763440860
$ws.Range("L16").Value = 8.080808080808
$ws.Range("M16").Value = -43.085106382978
$ws.Range("N16").Value = -82.136894824707
$ws.Range("C17").NumberFormat = "General"
$ws.Range("C17").Value = "'0"
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = -100
$ws.Range("G17").Value = 7
$ws.Range("H17").Value = 57.142857142857
$ws.Range("J17").Value = 140
$ws.Range("K17").Value = 10
$ws.Range("L17").Value = 4.761904761904
$ws.Range("M17").Value = 37.5
$ws.Range("N17").Value = -45.195729537366
$ws.Range("C18").Value = 2
$ws.Range("F18").Value = 8
$ws.Range("H18").Value = 60
$ws.Range("I18").Value = 65
$ws.Range("K18").Value = -17.721518987341
$ws.Range("L18").Value = -34.343434343434
$ws.Range("M18").Value = -73.360655737704
$ws.Range("N18").Value = -93.454179254783
$ws.Range("C19").Value = 12
$ws.Range("D19").Value = 12
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 40
$ws.Range("G19").Value = 54
$ws.Range("H19").Value = -25.925925925925
$ws.Range("I19").Value = 501
$ws.Range("J19").Value = 572
$ws.Range("K19").Value = -12.412587412587
$ws.Range("L19").Value = 12.080536912751
$ws.Range("M19").Value = 15.437788018433
$ws.Range("N19").Value = -10.535714285714
$ws.Range("C20").Value = 5
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = 150
$ws.Range("G20").Value = 20
$ws.Range("H20").Value = -30
$ws.Range("I20").Value = 155
$ws.Range("J20").Value = 117
$ws.Range("K20").Value = 32.478632478632
$ws.Range("L20").Value = 37.168141592920
$ws.Range("M20").Value = 9.154929577464
$ws.Range("N20").Value = -93.520066889632
$ws.Range("C21").Value = 26
$ws.Range("D21").Value = 18
$ws.Range("E21").Value = 44.444444444444
$ws.Range("F21").Value = 84
$ws.Range("G21").Value = 90
$ws.Range("H21").Value = -6.666666666666
$ws.Range("I21").Value = 992
$ws.Range("J21").Value = 1017
$ws.Range("K21").Value = -2.458210422812
$ws.Range("L21").Value = 8.533916849015
$ws.Range("M21").Value = -12.829525483304
$ws.Range("N21").Value = -79.592676404032
$ws.Range("C23").NumberFormat = "#,##0"
$ws.Range("C23").Value = 2
$ws.Range("D23").NumberFormat = "General"
$ws.Range("D23").Value = "'0"
$ws.Range("E23").NumberFormat = "General"
$ws.Range("E23").Value = "'***.*"
$ws.Range("F23").Value = 4
$ws.Range("H23").Value = -20
$ws.Range("I23").Value = 27
$ws.Range("K23").Value = 3.846153846153
$ws.Range("L23").Value = -6.896551724137
$ws.Range("M23").Value = 8
$ws.Range("C24").Value = 29
$ws.Range("D24").Value = 35
$ws.Range("E24").Value = -17.142857142857
$ws.Range("F24").Value = 115
$ws.Range("G24").Value = 114
$ws.Range("H24").Value = 0.877192982456
$ws.Range("I24").Value = 1412
$ws.Range("J24").Value = 1130
$ws.Range("K24").Value = 24.955752212389
$ws.Range("L24").Value = 42.770475227502
$ws.Range("M24").Value = 62.298850574712
$ws.Range("C25").Value = 23
$ws.Range("D25").Value = 29
$ws.Range("E25").Value = -20.689655172413
$ws.Range("F25").Value = 88
$ws.Range("G25").Value = 102
$ws.Range("H25").Value = -13.725490196078
$ws.Range("I25").Value = 1169
$ws.Range("J25").Value = 831
$ws.Range("K25").Value = 40.673886883273
$ws.Range("L25").Value = 86.443381180223
$ws.Range("C26").Value = 2
$ws.Range("D26").Value = 5
$ws.Range("E26").Value = -60
$ws.Range("F26").Value = 15
$ws.Range("G26").Value = 24
$ws.Range("H26").Value = -37.5
$ws.Range("I26").Value = 265
$ws.Range("J26").Value = 233
$ws.Range("K26").Value = 13.733905579399
$ws.Range("L26").Value = 29.901960784313
$ws.Range("M26").Value = -14.516129032258
$ws.Range("G27").Value = 3
$ws.Range("J27").Value = 23
$ws.Range("K27").Value = -39.130434782608
$ws.Range("L27").Value = 0
$ws.Range("C28").NumberFormat = "General"
$ws.Range("C28").Value = "'0"
$ws.Range("D28").Value = 2
$ws.Range("E28").Value = -100
$ws.Range("G28").Value = 4
$ws.Range("H28").Value = -25
$ws.Range("J28").Value = 27
$ws.Range("K28").Value = 33.333333333333
$ws.Range("L28").Value = 24.137931034482
$ws.Range("M29").Value = -83.333333333333
$ws.Range("M30").Value = -80
